$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$durations = @(513,763,676,307,320,282,592,873,325,343,772,720,358,640,794,644,372,1020,358,900,1126,320,338,407,651,1099,415,475,475,1018,786,1161,839,458,915,786,760,459,856,437,804,530,931,414,736,408,697,277,1134,762,697,489,415,224,746,577,455,550,486,319)

$ws.Range("P1").Value = "duration"
$ws.Range("P1").Style = $ws.Range("O1").Style

for ($i = 0; $i -lt $durations.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 16).Value = $durations[$i]
}
